# Apply the "Trade #13 closed" update to the live trading results workbook.
#
# Summary of changes:
#   - Summary sheet: update Current Capital, Total P&L $, Total P&L %,
#     Total Trades, Losing Trades, Win Rate % to reflect the new trade.
#   - Strategy Status sheet: update the MarketMaking row's Capital, Trades,
#     P&L $, P&L %, Win Rate % to reflect the new trade.
#   - All Trades sheet: append a new row (#14 / Trade 13) describing the
#     newly closed trade.
#   - MarketMaking sheet: append the same new row (#14 / Trade 13).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.03   # Current Capital
$summary.Range("B4").Value = 0.03      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 13        # Total Trades
$summary.Range("B8").Value = 5         # Losing Trades
$summary.Range("B9").Value = 46.15     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row = row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.03     # Capital
$status.Range("D4").Value = 13         # Trades
$status.Range("E4").Value = 0.03       # P&L $
$status.Range("F4").Value = 0.03       # P&L %
$status.Range("G4").Value = 46.15      # Win Rate %

# ---------------------------------------------------------------------
# Helper: append the new trade-13 row (row 14) to a trades-log sheet.
# Column B holds a literal date-like string ("2026-02-17"); Excel would
# normally auto-convert such text to a date serial number, so the cell's
# number format is forced to Text ("@") before the value is written so
# it remains a plain string, matching the rest of the sheet.
# ---------------------------------------------------------------------
function Add-Trade13Row($ws) {
    $ws.Cells.Item(14, 1).Value = 13
    $ws.Cells.Item(14, 2).NumberFormat = "@"
    $ws.Cells.Item(14, 2).Value = "2026-02-17"
    $ws.Cells.Item(14, 3).Value = "12:27:53"
    $ws.Cells.Item(14, 4).Value = "MarketMaking"
    $ws.Cells.Item(14, 5).Value = "UP"
    $ws.Cells.Item(14, 6).Value = 0.09
    $ws.Cells.Item(14, 7).Value = 0.07722
    $ws.Cells.Item(14, 8).Value = "CLOSED"
    $ws.Cells.Item(14, 9).Value = -14.1999
    $ws.Cells.Item(14, 10).Value = -0.01
    $ws.Cells.Item(14, 11).Value = 100.03
    $ws.Cells.Item(14, 12).Value = 0
    $ws.Cells.Item(14, 13).Value = 0
    $ws.Cells.Item(14, 14).Value = 0.6
    $ws.Cells.Item(14, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(14, 16).Value = "early_exit"
    $ws.Cells.Item(14, 17).Value = 0.14
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade13Row $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade13Row $marketMaking

Write-Host "Applied trade #13 closing update."
